# Generate Report for Handoff
# Updates the handoff identifiers (GUID changed from 74a5f089-... to
# c2877b6b-...), the corresponding target-file hashes, and the handoff
# timestamps across the Overview / zh-cn / de-de worksheets, keeping the
# external hyperlink addresses intact but refreshing their display text.

$wb = $excel.ActiveWorkbook

$oldGuid = "74a5f089-0559-4763-99d8-c3b941894a9a"
$newGuid = "c2877b6b-57f3-42d9-be73-981afda8b339"

$oldHash = "6317f33b8f86c1f97106924884dff4a7bcac1269"
$newHash = "be78c7147117786a0304fdcb816bd207132b18ec"

$newMdName = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newHandoffDate = "2016-03-22 12:32:44"
$newZhHandoffDatetime = "2016-03-22 12:32:34"
$newDeHandoffDatetime = "2016-03-22 12:32:44"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Address -like "*$oldGuid.md") {
        $hl.TextToDisplay = $newMdName
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhXlfName
$wsZhCn.Range("E2").Value = $newZhHandoffDatetime

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Address -like "*$oldGuid.md") {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.Address -like "*.zh-cn.xlf") {
        $hl.TextToDisplay = $newZhXlfName
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeXlfName
$wsDeDe.Range("E2").Value = $newDeHandoffDatetime

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Address -like "*$oldGuid.md") {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.Address -like "*.de-de.xlf") {
        $hl.TextToDisplay = $newDeXlfName
    }
}
